# cellExists() / getCell() should support UTF-8 named cells.
# Add a new greeting cell on Sheet1 (B4) holding a UTF-8 string, centered
# like the rest of the data block, and define a UTF-8 named range that
# refers to it.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# New cell with a UTF-8 (Devanagari) value, matching the existing
# centered style used by the rest of the A1:C3 block.
$cell = $ws1.Range("B4")
$cell.Value = "नमस्ते"
$cell.HorizontalAlignment = -4108
$cell.Select()

# Define a UTF-8 (Greek) named range pointing at the new cell. Add it
# with a plain ASCII name first, then rename it to the UTF-8 name.
$wb.Names.Add("TempUtf8Name", "=Sheet1!`$B`$4")
$wb.Names.Item("TempUtf8Name").Name = "Χαιρετισμός"
